$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 45.8803895
$ws.Range("H2").Value = 91.760779
$ws.Range("I2").Value = 0.04823508666161738
$ws.Range("J2").Value = 0.03353989414786691
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3657716666666667
$ws.Range("N2").Value = 1.097315
$ws.Range("O2").Value = 0.03059585711603819
$ws.Range("P2").Value = 0.03059585711603819
$ws.Range("Q2").Value = 16.78174653473083
$ws.Range("R2").Value = 100.690479208385
$ws.Range("S2").Value = 0.001475793819478565
$ws.Range("T2").Value = 0.001026181809035181

$ws.Range("G3").Value = 45.8803895
$ws.Range("H3").Value = 91.760779
$ws.Range("I3").Value = 0.04823508666161738
$ws.Range("J3").Value = 0.03353989414786691
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 11.510206
$ws.Range("N3").Value = 34.530618
$ws.Range("O3").Value = 0.9627990635838353
$ws.Range("P3").Value = 0.9627990635838353
$ws.Range("Q3").Value = 528.0927345052369
$ws.Range("R3").Value = 3168.556407031422
$ws.Range("S3").Value = 0.04644069626969035
$ws.Range("T3").Value = 0.03229217867826722

$ws.Range("G4").Value = 45.8803895
$ws.Range("H4").Value = 91.760779
$ws.Range("I4").Value = 0.04823508666161738
$ws.Range("J4").Value = 0.03353989414786691
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.07896333333333333
$ws.Range("N4").Value = 0.23689
$ws.Range("O4").Value = 0.006605079300126477
$ws.Range("P4").Value = 0.006605079300126477
$ws.Range("Q4").Value = 3.622868489551666
$ws.Range("R4").Value = 21.73721093731
$ws.Range("S4").Value = 0.0003185965724484557
$ws.Range("T4").Value = 0.0002215336605645089

$ws.Range("I5").Value = 0.5969743022710003
$ws.Range("J5").Value = 0.6226521903328125
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3657716666666667
$ws.Range("N5").Value = 1.097315
$ws.Range("O5").Value = 0.03059585711603819
$ws.Range("P5").Value = 0.03059585711603819
$ws.Range("Q5").Value = 207.6967643645111
$ws.Range("R5").Value = 1869.2708792806
$ws.Range("S5").Value = 0.01826494045423012
$ws.Range("T5").Value = 0.01905057744841094

$ws.Range("I6").Value = 0.5969743022710003
$ws.Range("J6").Value = 0.6226521903328125
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 11.510206
$ws.Range("N6").Value = 34.530618
$ws.Range("O6").Value = 0.9627990635838353
$ws.Range("P6").Value = 0.9627990635838353
$ws.Range("Q6").Value = 6535.860377473145
$ws.Range("R6").Value = 58822.74339725832
$ws.Range("S6").Value = 0.5747662992101326
$ws.Range("T6").Value = 0.5994889457908559

$ws.Range("I7").Value = 0.5969743022710003
$ws.Range("J7").Value = 0.6226521903328125
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.07896333333333333
$ws.Range("N7").Value = 0.23689
$ws.Range("O7").Value = 0.006605079300126477
$ws.Range("P7").Value = 0.006605079300126477
$ws.Range("Q7").Value = 44.83788748928888
$ws.Range("R7").Value = 403.5409874036
$ws.Range("S7").Value = 0.003943062606637631
$ws.Range("T7").Value = 0.004112667093545672

$ws.Range("G8").Value = 58.29942466666666
$ws.Range("H8").Value = 174.898274
$ws.Range("I8").Value = 0.06129149799652638
$ws.Range("J8").Value = 0.0639278530602342
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.3657716666666667
$ws.Range("N8").Value = 1.097315
$ws.Range("O8").Value = 0.03059585711603819
$ws.Range("P8").Value = 0.03059585711603819
$ws.Range("Q8").Value = 21.32427772603444
$ws.Range("R8").Value = 191.91849953431
$ws.Range("S8").Value = 0.001875265915129662
$ws.Range("T8").Value = 0.00195592745796601

$ws.Range("G9").Value = 58.29942466666666
$ws.Range("H9").Value = 174.898274
$ws.Range("I9").Value = 0.06129149799652638
$ws.Range("J9").Value = 0.0639278530602342
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.510206
$ws.Range("N9").Value = 34.530618
$ws.Range("O9").Value = 0.9627990635838353
$ws.Range("P9").Value = 0.9627990635838353
$ws.Range("Q9").Value = 671.0383875948145
$ws.Range("R9").Value = 6039.345488353331
$ws.Range("S9").Value = 0.05901139687670612
$ws.Range("T9").Value = 0.06154967706331851

$ws.Range("G10").Value = 58.29942466666666
$ws.Range("H10").Value = 174.898274
$ws.Range("I10").Value = 0.06129149799652638
$ws.Range("J10").Value = 0.0639278530602342
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.07896333333333333
$ws.Range("N10").Value = 0.23689
$ws.Range("O10").Value = 0.006605079300126477
$ws.Range("P10").Value = 0.006605079300126477
$ws.Range("Q10").Value = 4.603516903095555
$ws.Range("R10").Value = 41.43165212786
$ws.Range("S10").Value = 0.0004048352046905998
$ws.Range("T10").Value = 0.00042224853894968

$ws.Range("G11").Value = 71.7986375
$ws.Range("H11").Value = 143.597275
$ws.Range("I11").Value = 0.07548352443691768
$ws.Range("J11").Value = 0.05248688443918001
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.3657716666666667
$ws.Range("N11").Value = 1.097315
$ws.Range("O11").Value = 0.03059585711603819
$ws.Range("P11").Value = 0.03059585711603819
$ws.Range("Q11").Value = 26.26190730277083
$ws.Range("R11").Value = 157.571443816625
$ws.Range("S11").Value = 0.00230948312828691
$ws.Range("T11").Value = 0.00160588121676716

$ws.Range("G12").Value = 71.7986375
$ws.Range("H12").Value = 143.597275
$ws.Range("I12").Value = 0.07548352443691768
$ws.Range("J12").Value = 0.05248688443918001
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 11.510206
$ws.Range("N12").Value = 34.530618
$ws.Range("O12").Value = 0.9627990635838353
$ws.Range("P12").Value = 0.9627990635838353
$ws.Range("Q12").Value = 826.4171081443249
$ws.Range("R12").Value = 4958.50264886595
$ws.Range("S12").Value = 0.07267546664387189
$ws.Range("T12").Value = 0.05053432318847549

$ws.Range("G13").Value = 71.7986375
$ws.Range("H13").Value = 143.597275
$ws.Range("I13").Value = 0.07548352443691768
$ws.Range("J13").Value = 0.05248688443918001
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.07896333333333333
$ws.Range("N13").Value = 0.23689
$ws.Range("O13").Value = 0.006605079300126477
$ws.Range("P13").Value = 0.006605079300126477
$ws.Range("Q13").Value = 5.669459745791666
$ws.Range("R13").Value = 34.01675847475
$ws.Range("S13").Value = 0.0004985746647588761
$ws.Range("T13").Value = 0.0003466800339373584

$ws.Range("G14").Value = 147.7885923333333
$ws.Range("H14").Value = 443.365777
$ws.Range("I14").Value = 0.1553734751706233
$ws.Range("J14").Value = 0.1620566149440249
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3657716666666667
$ws.Range("N14").Value = 1.097315
$ws.Range("O14").Value = 0.03059585711603819
$ws.Range("P14").Value = 0.03059585711603819
$ws.Range("Q14").Value = 54.05687973208389
$ws.Range("R14").Value = 486.511917588755
$ws.Range("S14").Value = 0.004753784645942697
$ws.Range("T14").Value = 0.004958261035536204

$ws.Range("G15").Value = 147.7885923333333
$ws.Range("H15").Value = 443.365777
$ws.Range("I15").Value = 0.1553734751706233
$ws.Range("J15").Value = 0.1620566149440249
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 11.510206
$ws.Range("N15").Value = 34.530618
$ws.Range("O15").Value = 0.9627990635838353
$ws.Range("P15").Value = 0.9627990635838353
$ws.Range("Q15").Value = 1701.077142206687
$ws.Range("R15").Value = 15309.69427986018
$ws.Range("S15").Value = 0.1495934364000424
$ws.Range("T15").Value = 0.1560279571156733

$ws.Range("G16").Value = 147.7885923333333
$ws.Range("H16").Value = 443.365777
$ws.Range("I16").Value = 0.1553734751706233
$ws.Range("J16").Value = 0.1620566149440249
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.07896333333333333
$ws.Range("N16").Value = 0.23689
$ws.Range("O16").Value = 0.006605079300126477
$ws.Range("P16").Value = 0.006605079300126477
$ws.Range("Q16").Value = 11.66987987928111
$ws.Range("R16").Value = 105.02891891353
$ws.Range("S16").Value = 0.001026254124638199
$ws.Range("T16").Value = 0.001070396792815346

$ws.Range("G17").Value = 59.584107
$ws.Range("H17").Value = 178.752321
$ws.Range("I17").Value = 0.06264211346331491
$ws.Range("J17").Value = 0.06533656307588157
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3657716666666667
$ws.Range("N17").Value = 1.097315
$ws.Range("O17").Value = 0.03059585711603819
$ws.Range("P17").Value = 0.03059585711603819
$ws.Range("Q17").Value = 21.794178124235
$ws.Range("R17").Value = 196.147603118115
$ws.Range("S17").Value = 0.001916589152970235
$ws.Range("T17").Value = 0.001999028148322689

$ws.Range("G18").Value = 59.584107
$ws.Range("H18").Value = 178.752321
$ws.Range("I18").Value = 0.06264211346331491
$ws.Range("J18").Value = 0.06533656307588157
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 11.510206
$ws.Range("N18").Value = 34.530618
$ws.Range("O18").Value = 0.9627990635838353
$ws.Range("P18").Value = 0.9627990635838353
$ws.Range("Q18").Value = 685.8253458960419
$ws.Range("R18").Value = 6172.428113064378
$ws.Range("S18").Value = 0.06031176818339196
$ws.Range("T18").Value = 0.06290598174724496

$ws.Range("G19").Value = 59.584107
$ws.Range("H19").Value = 178.752321
$ws.Range("I19").Value = 0.06264211346331491
$ws.Range("J19").Value = 0.06533656307588157
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.07896333333333333
$ws.Range("N19").Value = 0.23689
$ws.Range("O19").Value = 0.006605079300126477
$ws.Range("P19").Value = 0.006605079300126477
$ws.Range("Q19").Value = 4.704959702409999
$ws.Range("R19").Value = 42.34463732168999
$ws.Range("S19").Value = 0.0004137561269527155
$ws.Range("T19").Value = 0.0004315531803139133
